# Update the "Marking"/"Total" row values on the marksheet's "quiz" sheet
# to reflect the corrected total correct/total marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 45
$ws.Range("E12").Value = "45/140"
